$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a numeric-looking string value must be forced to Text format
# first, otherwise Excel auto-converts them to numbers and the original
# formatting (trailing zeros, etc.) is lost.

# Row 2
$ws.Range("D2").Value = "70.467.79"
$ws.Range("E2").Value = "  -2.94%  "

# Row 3
$ws.Range("D3").Value = "3.917.87"
$ws.Range("E3").Value = "  -3.11%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.51"
$ws.Range("E5").Value = "  +3.87%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.94"
$ws.Range("E6").Value = "  -0.13%  "

# Row 7
$ws.Range("E7").Value = "  -5.05%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.732"
$ws.Range("E9").Value = "  -4.02%  "

# Row 10
$ws.Range("E10").Value = "  -4.98%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.70"
$ws.Range("E11").Value = "  +13.20%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000315"
$ws.Range("E12").Value = "  -2.84%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.46"
$ws.Range("E13").Value = "  -3.60%  "

# Row 14
$ws.Range("D14").Value = "4.540.60"
$ws.Range("E14").Value = "  -3.02%  "

# Row 15
$ws.Range("D15").Value = "3.922.40"
$ws.Range("E15").Value = "  -2.99%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.85"
$ws.Range("E16").Value = "  -2.04%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.18"
$ws.Range("E17").Value = "  -4.66%  "

# Row 18
$ws.Range("E18").Value = "  -0.77%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.16"
$ws.Range("E19").Value = "  -3.92%  "

# Row 20
$ws.Range("D20").Value = "70.343.30"
$ws.Range("E20").Value = "  -2.72%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "426.68"
$ws.Range("E21").Value = "  -3.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "96.42"
$ws.Range("E22").Value = "  -7.77%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.49"
$ws.Range("E23").Value = "  -1.72%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.17"
$ws.Range("E24").Value = "  +5.23%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.12"
$ws.Range("E25").Value = "  -3.80%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.13"
$ws.Range("E26").Value = "  -3.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.47"
$ws.Range("E27").Value = "  -5.95%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.85"
$ws.Range("E28").Value = "  +0.76%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.63"
$ws.Range("E29").Value = "  +16.90%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.13"
$ws.Range("E30").Value = "  -5.77%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.40"
$ws.Range("E31").Value = "  +9.10%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.31"
$ws.Range("E32").Value = "  -3.02%  "

# Row 33
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "676.24"
$ws.Range("E33").Value = "  -0.41%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.128"
$ws.Range("E34").Value = "  +0.30%  "

# Row 35
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "47.63"
$ws.Range("E35").Value = "  +14.82%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.03"
$ws.Range("E36").Value = "  -3.48%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.426"
$ws.Range("E37").Value = "  -0.76%  "

# Row 38
$ws.Range("E38").Value = "  -5.53%  "

# Row 39
$ws.Range("B39").Value = "ThetaToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.40"
$ws.Range("E39").Value = "  -3.75%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.147"
$ws.Range("E40").Value = "  -1.75%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.06%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.31"
$ws.Range("E42").Value = "  +4.06%  "

# Row 43
$ws.Range("E43").Value = "  +0.44%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0478"
$ws.Range("E44").Value = "  -2.80%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.147"
$ws.Range("E45").Value = "  -5.77%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.67"
$ws.Range("E46").Value = "  -2.31%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.52"
$ws.Range("E47").Value = "  +3.53%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.36"
$ws.Range("E48").Value = "  -3.76%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.95"
$ws.Range("E49").Value = "  -4.07%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000270"
$ws.Range("E50").Value = "  +1.51%  "

# Row 51
$ws.Range("E51").Value = "  +1.23%  "
